$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula-driven rows 10 and 11 with their computed static values
# (emulates a Copy > Paste Special > Values operation performed in the source edit)

$ws.Range("B10").Value = 366444211.24861097
$ws.Range("C10").Value = 75625581.995992005
$ws.Range("D10").Value = 71004194.989956662
$ws.Range("E10").Value = 186359759.21892264
$ws.Range("F10").Value = 62465454.06821765
$ws.Range("G10").Value = 1230140281.2897937
$ws.Range("H10").Value = 357096623.61340791
$ws.Range("I10").Value = 90753477.777032003
$ws.Range("J10").Value = 15029260.617756
$ws.Range("K10").Value = 2765801.9299363601
$ws.Range("L10").Value = 6934465.5969321597
$ws.Range("M10").Value = 121801119.7484466

$ws.Range("B11").Value = 51025890.865032718
$ws.Range("C11").Value = 6079079.2382131144
$ws.Range("D11").Value = 26611632.772054341
$ws.Range("E11").Value = 438544.10012470215
$ws.Range("F11").Value = 7562569.7118316786
$ws.Range("G11").Value = 26271025.787142672
$ws.Range("H11").Value = 8916885.493132893
$ws.Range("I11").Value = 2145244.1245425958
$ws.Range("J11").Value = 3033747.90745488
$ws.Range("K11").Value = 859565.21714523504
$ws.Range("L11").Value = 5136840.3628080003
$ws.Range("M11").Value = 6498745.6029703859

# Update the active cell selection shown in the saved sheet view
$ws.Range("N19").Select()
